$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 261.380203
$ws.Range("N2").Value = 784.1406089999999
$ws.Range("O2").Value = 0.6968677182772199
$ws.Range("P2").Value = 0.6968677182772199
$ws.Range("Q2").Value = 1587.134833422593
$ws.Range("R2").Value = 14284.21350080334
$ws.Range("S2").Value = 0.002748445770106722
$ws.Range("T2").Value = 0.002748445770106721

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.999428
$ws.Range("N3").Value = 95.998284
$ws.Range("O3").Value = 0.08531391482826334
$ws.Range("P3").Value = 0.08531391482826335
$ws.Range("Q3").Value = 194.304718741068
$ws.Range("R3").Value = 1748.742468669612
$ws.Range("S3").Value = 0.0003364780175506811
$ws.Range("T3").Value = 0.0003364780175506811

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.072131
$ws.Range("H4").Value = 18.216393
$ws.Range("I4").Value = 0.003943999267036455
$ws.Range("J4").Value = 0.003943999267036454
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 81.69901900000001
$ws.Range("N4").Value = 245.097057
$ws.Range("O4").Value = 0.2178183668945166
$ws.Range("P4").Value = 0.2178183668945167
$ws.Range("Q4").Value = 496.087145939489
$ws.Range("R4").Value = 4464.784313455401
$ws.Range("S4").Value = 0.0008590754793790512
$ws.Range("T4").Value = 0.0008590754793790512

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 261.380203
$ws.Range("N5").Value = 784.1406089999999
$ws.Range("O5").Value = 0.6968677182772199
$ws.Range("P5").Value = 0.6968677182772199
$ws.Range("Q5").Value = 387065.3458394501
$ws.Range("R5").Value = 3483588.11255505
$ws.Range("S5").Value = 0.6702821273434149
$ws.Range("T5").Value = 0.6702821273434147

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1480.851806666667
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9618498744646554
$ws.Range("J6").Value = 0.9618498744646552
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 31.999428
$ws.Range("N6").Value = 95.998284
$ws.Range("O6").Value = 0.08531391482826334
$ws.Range("P6").Value = 0.08531391482826335
$ws.Range("Q6").Value = 47386.41076609991
$ws.Range("R6").Value = 426477.6968948992
$ws.Range("S6").Value = 0.0820591782676534
$ws.Range("T6").Value = 0.0820591782676534

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1480.851806666667
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9618498744646554
$ws.Range("J7").Value = 0.9618498744646552
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 81.69901900000001
$ws.Range("N7").Value = 245.097057
$ws.Range("O7").Value = 0.2178183668945166
$ws.Range("P7").Value = 0.2178183668945167
$ws.Range("Q7").Value = 120984.1398890443
$ws.Range("R7").Value = 1088857.259001399
$ws.Range("S7").Value = 0.2095085688535871
$ws.Range("T7").Value = 0.2095085688535871

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 52.663316
$ws.Range("H8").Value = 157.989948
$ws.Range("I8").Value = 0.03420612626830831
$ws.Range("J8").Value = 0.0342061262683083
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 261.380203
$ws.Range("N8").Value = 784.1406089999999
$ws.Range("O8").Value = 0.6968677182772199
$ws.Range("P8").Value = 0.6968677182772199
$ws.Range("Q8").Value = 13765.14822673315
$ws.Range("R8").Value = 123886.3340405983
$ws.Range("S8").Value = 0.02383714516369849
$ws.Range("T8").Value = 0.02383714516369848

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 52.663316
$ws.Range("H9").Value = 157.989948
$ws.Range("I9").Value = 0.03420612626830831
$ws.Range("J9").Value = 0.0342061262683083
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 31.999428
$ws.Range("N9").Value = 95.998284
$ws.Range("O9").Value = 0.08531391482826334
$ws.Range("P9").Value = 0.08531391482826335
$ws.Range("Q9").Value = 1685.195988583248
$ws.Range("R9").Value = 15166.76389724923
$ws.Range("S9").Value = 0.002918258543059276
$ws.Range("T9").Value = 0.002918258543059276

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 52.663316
$ws.Range("H10").Value = 157.989948
$ws.Range("I10").Value = 0.03420612626830831
$ws.Range("J10").Value = 0.0342061262683083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 81.69901900000001
$ws.Range("N10").Value = 245.097057
$ws.Range("O10").Value = 0.2178183668945166
$ws.Range("P10").Value = 0.2178183668945167
$ws.Range("Q10").Value = 4302.541254487004
$ws.Range("R10").Value = 38722.87129038304
$ws.Range("S10").Value = 0.007450722561550543
$ws.Range("T10").Value = 0.007450722561550542
